# Teacher availability constraint added
# Swap "Dig. Signal" (2:00-3:00 slot, row 22) and "Dig. Analytics" (3:00-4:00 slot, row 25)
# and rename "R & Python Lab" -> "R & Python" (row 25 and the course list row 33),
# plus drop the trailing "Lab" line from the course-name cell in row 33.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 22 (2:00 PM - 3:00 PM slot): Dig. Signal -> Dig. Analytics
$ws.Range("D22").Value = "Dig. Analytics"
$ws.Range("E22").Value = "Dig. Analytics"
$ws.Range("F22").Value = "Dig. Analytics"

# Row 25 (3:00 PM - 4:00 PM slot): Dig. Analytics -> Dig. Signal
$ws.Range("D25").Value = "Dig. Signal"
$ws.Range("E25").Value = "Dig. Signal"
$ws.Range("F25").Value = "Dig. Signal"

# Row 25: course-id label rename
$ws.Range("B25").Value = "R & Python"

# Row 33: course-id label rename + course-name cell no longer has a "Lab" second line
$ws.Range("B33").Value = "R & Python"
$ws.Range("C33").Value = "R and Python Programming"
